# This workbook's data rows (2-31) were reshuffled by an upstream weekly
# data refresh: every row's full record (Fecha, Variedad, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg)
# moved to a different row position while columns that don't vary
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Calidad,
# Kg o Unidades, Clasificacion) stayed constant across every row anyway.
#
# targetRow -> sourceRow (source row's CURRENT/"before" content becomes
# the target row's new content)
$rowMap = @{
    2  = 14
    3  = 21
    4  = 26
    5  = 10
    6  = 13
    7  = 30
    8  = 23
    9  = 27
    10 = 19
    11 = 7
    12 = 18
    13 = 3
    14 = 6
    15 = 9
    16 = 8
    17 = 2
    18 = 16
    19 = 24
    20 = 25
    21 = 4
    22 = 15
    23 = 5
    24 = 29
    25 = 20
    26 = 12
    27 = 31
    28 = 17
    29 = 22
    30 = 11
    31 = 28
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ row-to-row (and therefore need to be
# carried along with the permutation). Columns A,B,C,E,F,G,I,Q,R are
# constant across all data rows, so they don't need to move.
$numericCols = @("D", "J", "K", "L", "M", "P")
$textCols = @("H", "N", "O")

# 1) Snapshot every data row's current values BEFORE writing anything,
#    since several rows both give and receive data (a true permutation).
$snapshot = @{}
for ($r = 2; $r -le 31; $r++) {
    $row = @{}
    foreach ($col in $numericCols) {
        $row[$col] = $ws.Range("$col$r").Value2
    }
    foreach ($col in $textCols) {
        $row[$col] = $ws.Range("$col$r").Text
    }
    $snapshot[$r] = $row
}

# 2) Write each target row's new content from its mapped source row's
#    snapshot.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $srcData = $snapshot[$sourceRow]
    foreach ($col in $numericCols) {
        $ws.Range("$col$targetRow").Value = $srcData[$col]
    }
    foreach ($col in $textCols) {
        $ws.Range("$col$targetRow").Value = $srcData[$col]
    }
}
